# Insert a new weekly data row at the top of the Mango price list
# (row 134), pushing the existing rows 134:230 down to 135:231.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 134, shifting rows 134:230 down one row.
$ws.Rows("134:134").Insert()

# Populate the newly inserted row 134 with the new weekly record.
$ws.Range("A134").Value = 4
$ws.Range("B134").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C134").Value = "Los Lagos"
$ws.Range("D134").Value = 44790
$ws.Range("D134").NumberFormat = $ws.Range("D135").NumberFormat
$ws.Range("E134").Value = 10
$ws.Range("F134").Value = "Fruta"
$ws.Range("G134").Value = 100108
$ws.Range("H134").Value = "Tropicales y subtropicales"
$ws.Range("I134").Value = 100108002
$ws.Range("J134").Value = "Mango"
$ws.Range("K134").Value = "Sin especificar"
$ws.Range("L134").Value = "Primera"
$ws.Range("M134").Value = 60
$ws.Range("N134").Value = 13000
$ws.Range("O134").Value = 14000
$ws.Range("P134").Value = 13500
$ws.Range("Q134").Value = "`$/bandeja 4 kilos"
$ws.Range("R134").Value = "Brasil"
$ws.Range("S134").Value = 3375
$ws.Range("T134").Value = 4
